$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 11:57:36"
$ws1.Cells.Item(3, 1).Value = "Total filas: 202"

$ws1.Cells.Item(189, 2).Value = "11:57:25"
$ws1.Cells.Item(189, 3).Value = "12:02"
$ws1.Cells.Item(189, 4).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(189, 5).Value = 5
$ws1.Cells.Item(189, 6).Value = "LP1912"
$ws1.Cells.Item(189, 7).Value = "30/12/2025"

$ws1.Cells.Item(190, 2).Value = "11:57:25"
$ws1.Cells.Item(190, 3).Value = "12:07"
$ws1.Cells.Item(190, 4).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(190, 5).Value = 10
$ws1.Cells.Item(190, 6).Value = "LP1912"
$ws1.Cells.Item(190, 7).Value = "30/12/2025"

$ws1.Cells.Item(191, 2).Value = "11:57:25"
$ws1.Cells.Item(191, 3).Value = "12:07"
$ws1.Cells.Item(191, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(191, 5).Value = 10
$ws1.Cells.Item(191, 6).Value = "LP1912"
$ws1.Cells.Item(191, 7).Value = "30/12/2025"

$ws1.Cells.Item(192, 2).Value = "11:57:25"
$ws1.Cells.Item(192, 3).Value = "12:21"
$ws1.Cells.Item(192, 4).Value = "14_ABASTO"
$ws1.Cells.Item(192, 5).Value = 24
$ws1.Cells.Item(192, 6).Value = "LP1912"
$ws1.Cells.Item(192, 7).Value = "30/12/2025"

$ws1.Cells.Item(193, 2).Value = "11:57:25"
$ws1.Cells.Item(193, 3).Value = "12:21"
$ws1.Cells.Item(193, 4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(193, 5).Value = 24
$ws1.Cells.Item(193, 6).Value = "LP1912"
$ws1.Cells.Item(193, 7).Value = "30/12/2025"

$ws1.Cells.Item(194, 2).Value = "11:57:25"
$ws1.Cells.Item(194, 3).Value = "12:35"
$ws1.Cells.Item(194, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(194, 5).Value = 38
$ws1.Cells.Item(194, 6).Value = "LP1912"
$ws1.Cells.Item(194, 7).Value = "30/12/2025"

$ws1.Cells.Item(195, 2).Value = "11:57:25"
$ws1.Cells.Item(195, 3).Value = "12:38"
$ws1.Cells.Item(195, 4).Value = "17_179 Y 38"
$ws1.Cells.Item(195, 5).Value = 41
$ws1.Cells.Item(195, 6).Value = "LP1912"
$ws1.Cells.Item(195, 7).Value = "30/12/2025"

$ws1.Cells.Item(196, 2).Value = "11:57:25"
$ws1.Cells.Item(196, 3).Value = "12:48"
$ws1.Cells.Item(196, 4).Value = "10_OLMOS"
$ws1.Cells.Item(196, 5).Value = 51
$ws1.Cells.Item(196, 6).Value = "LP1912"
$ws1.Cells.Item(196, 7).Value = "30/12/2025"

$ws1.Cells.Item(197, 2).Value = "11:57:25"
$ws1.Cells.Item(197, 3).Value = "12:49"
$ws1.Cells.Item(197, 4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(197, 5).Value = 52
$ws1.Cells.Item(197, 6).Value = "LP1912"
$ws1.Cells.Item(197, 7).Value = "30/12/2025"

$ws1.Cells.Item(198, 2).Value = "11:57:25"
$ws1.Cells.Item(198, 3).Value = "12:50"
$ws1.Cells.Item(198, 4).Value = "15_ABASTO"
$ws1.Cells.Item(198, 5).Value = 53
$ws1.Cells.Item(198, 6).Value = "LP1912"
$ws1.Cells.Item(198, 7).Value = "30/12/2025"

$ws1.Cells.Item(199, 2).Value = "11:57:25"
$ws1.Cells.Item(199, 3).Value = "12:55"
$ws1.Cells.Item(199, 4).Value = "10_OLMOS"
$ws1.Cells.Item(199, 5).Value = 58
$ws1.Cells.Item(199, 6).Value = "LP1912"
$ws1.Cells.Item(199, 7).Value = "30/12/2025"

$ws1.Cells.Item(200, 2).Value = "11:57:25"
$ws1.Cells.Item(200, 3).Value = "13:07"
$ws1.Cells.Item(200, 4).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(200, 5).Value = 70
$ws1.Cells.Item(200, 6).Value = "LP1912"
$ws1.Cells.Item(200, 7).Value = "30/12/2025"

$ws1.Cells.Item(201, 2).Value = "11:57:25"
$ws1.Cells.Item(201, 3).Value = "13:27"
$ws1.Cells.Item(201, 4).Value = "14_ABASTO"
$ws1.Cells.Item(201, 5).Value = 90
$ws1.Cells.Item(201, 6).Value = "LP1912"
$ws1.Cells.Item(201, 7).Value = "30/12/2025"

$ws1.Cells.Item(202, 2).Value = "11:57:25"
$ws1.Cells.Item(202, 3).Value = "13:27"
$ws1.Cells.Item(202, 4).Value = "17_ROMERO"
$ws1.Cells.Item(202, 5).Value = 90
$ws1.Cells.Item(202, 6).Value = "LP1912"
$ws1.Cells.Item(202, 7).Value = "30/12/2025"

$ws1.Cells.Item(203, 2).Value = "11:57:25"
$ws1.Cells.Item(203, 3).Value = "13:35"
$ws1.Cells.Item(203, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(203, 5).Value = 98
$ws1.Cells.Item(203, 6).Value = "LP1912"
$ws1.Cells.Item(203, 7).Value = "30/12/2025"

# ---- Sheet 2: LP1912-215 (timestamp only) ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 11:57:36"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 11:57:36"
$ws3.Cells.Item(3, 1).Value = "Total filas: 28"

$ws3.Cells.Item(27, 2).Value = "30/12/2025"
$ws3.Cells.Item(27, 3).Value = "11:57:36"
$ws3.Cells.Item(27, 4).Value = "12:04"
$ws3.Cells.Item(27, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(27, 6).Value = 7
$ws3.Cells.Item(27, 7).Value = "L6173"

$ws3.Cells.Item(28, 2).Value = "30/12/2025"
$ws3.Cells.Item(28, 3).Value = "11:57:31"
$ws3.Cells.Item(28, 4).Value = "12:54"
$ws3.Cells.Item(28, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(28, 6).Value = 57
$ws3.Cells.Item(28, 7).Value = "L6203"

$ws3.Cells.Item(29, 2).Value = "30/12/2025"
$ws3.Cells.Item(29, 3).Value = "11:57:36"
$ws3.Cells.Item(29, 4).Value = "13:31"
$ws3.Cells.Item(29, 5).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(29, 6).Value = 94
$ws3.Cells.Item(29, 7).Value = "L6173"

